$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.965.46"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "'2.548.14"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'569.18"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "'146.49"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.583"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "'2.545.77"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  -4.12%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "'27.19"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "'3.000.85"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "'62.883.92"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'2.536.60"
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("D19").Value = "'11.30"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'4.35"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "'334.72"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'65.37"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("E26").Value = "  +4.68%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("B28").Value = "SuiNetwork"
$ws.Range("C28").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D28").Value = "'1.47"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'8.34"
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").Value = "'7.35"
$ws.Range("E30").Value = "  +8.70%  "
$ws.Range("D31").Value = "'0.0₃0813"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "'175.50"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").Value = "'404.93"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("D36").Value = "'19.13"
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D42").Value = "'39.45"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").Value = "'151.82"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "'20.73"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "'0.0532"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("D47").Value = "'0.601"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.0966"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0240"
$ws.Range("E49").Value = "  +4.83%  "
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  -2.93%  "
